# This script rewrites the player table (A2:C19) on the active sheet so
# that it matches the updated "Team of Outs" roster: rows are reordered,
# "Dillon Brooks / SG,SF / Houston Rockets" is removed, and a new row for
# "Jonas Valanciunas / C / Washington Wizards" is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ja Morant",          "PG",          "Memphis Grizzlies"),
    @("Tyler Herro",        "PG,SG",       "Miami Heat"),
    @("Josh Giddey",        "PG,SG,SF",    "Chicago Bulls"),
    @("Scottie Barnes",     "PG,SG,SF,PF", "Toronto Raptors"),
    @("DeMar DeRozan",      "SF,PF",       "Sacramento Kings"),
    @("Miles Bridges",      "SF,PF",       "Charlotte Hornets"),
    @("Shaedon Sharpe",     "SG,SF",       "Portland Trail Blazers"),
    @("Nikola Vucevic",     "PF,C",        "Chicago Bulls"),
    @("Brook Lopez",        "C",           "Milwaukee Bucks"),
    @("Kevon Looney",       "PF,C",        "Golden State Warriors"),
    @("Jonas Valanciunas",  "C",           "Washington Wizards"),
    @("Nick Richards",      "C",           "Phoenix Suns"),
    @("De'Aaron Fox",       "PG",          "Sacramento Kings"),
    @("Bobby Portis",       "PF,C",        "Milwaukee Bucks"),
    @("Mikal Bridges",      "SG,SF,PF",    "New York Knicks"),
    @("Luka Doncic",        "PG,SG",       "Dallas Mavericks"),
    @("Evan Mobley",        "PF,C",        "Cleveland Cavaliers"),
    @("Clint Capela",       "C",           "Atlanta Hawks")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
